$d = $word.ActiveDocument

$replacements = @(
    @("778×3=2334", "886×4=3544"),
    @("406×9=3654", "595×8=4760"),
    @("518×9=4662", "448×7=3136"),
    @("567×7=3969", "361×6=2166"),
    @("166×9=1494", "980×5=4900"),
    @("859×6=5154", "175×5=875"),
    @("283×3=849", "304×4=1216"),
    @("338×7=2366", "551×9=4959"),
    @("567×3=1701", "705×9=6345"),
    @("556×3=1668", "551×7=3857"),
    @("965×8=7720", "718×6=4308"),
    @("247×3=741", "120×3=360"),
    @("947×4=3788", "405×2=810"),
    @("123×8=984", "286×7=2002"),
    @("410×3=1230", "692×2=1384"),
    @("643×6=3858", "214×5=1070"),
    @("460×6=2760", "776×4=3104"),
    @("101×7=707", "504×6=3024"),
    @("656×7=4592", "597×3=1791"),
    @("713×9=6417", "955×3=2865"),
    @("481×3=1443", "916×6=5496"),
    @("649×4=2596", "775×2=1550"),
    @("831×8=6648", "477×3=1431"),
    @("247×4=988", "916×2=1832"),
    @("850×6=5100", "225×3=675")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
